# EV-141 Add party abbr for seeding
# Adds a new "Skratka" (abbreviation) column to the PAR_2020_tab0a table,
# populating it with short party abbreviations for seeding purposes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# 1) Grow the table by one column (becomes column H, A1:H26).
$newCol = $tbl.ListColumns.Add()

# 2) Header - setting the cell value syncs the ListColumn/table name too.
$ws.Range("H1").Value = "Skratka"

# 3) Row 12 (OĽANO) already had a pre-existing "special" font style sitting
#    on the empty H12 cell (same family as the highlighted G12/G6/G25 rows).
#    Mirror column G's formatting for that row before writing its value so
#    the abbreviation keeps the same distinctive look as the Farba column.
$ws.Range("G12").Copy() | Out-Null
$ws.Range("H12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# 4) Fill in the abbreviations, row by row.
$ws.Range("H2").Value = "SĽS"
$ws.Range("H3").Value = "DOBRÁ VOĽBA"
$ws.Range("H4").Value = "SAS"
$ws.Range("H5").Value = "SME RODINA"
$ws.Range("H6").Value = "SHO"
$ws.Range("H7").Value = "ZA ĽUDÍ"
$ws.Range("H8").Value = "MÁME TOHO DOSŤ"
$ws.Range("H9").Value = "HLAS PRAV"
$ws.Range("H10").Value = "SNS"
$ws.Range("H11").Value = "DS"
$ws.Range("H12").Value = "OĽANO"
$ws.Range("H13").Value = "PS SPOLU"
$ws.Range("H14").Value = "STANK"

# Row 15 is seeded as a ratio (99%) rather than a text abbreviation.
$ws.Range("H15").Value = 0.99
$ws.Range("H15").NumberFormat = "0%"

$ws.Range("H16").Value = "KDH"
$ws.Range("H17").Value = "SLOV. LIGA"
$ws.Range("H18").Value = "VLASŤ"
$ws.Range("H19").Value = "MOST-HID"
$ws.Range("H20").Value = "SMER-SD"
$ws.Range("H21").Value = "SOLIDARITA"
$ws.Range("H22").Value = "HLAS ĽUDU"
$ws.Range("H23").Value = "MKO-MKS"
$ws.Range("H24").Value = "PSN"
$ws.Range("H25").Value = "ĽSNS"
$ws.Range("H26").Value = "SOCIALISTI"

# 5) Column H sizing - roomy enough for the longest abbreviation.
$ws.Columns.Item(8).ColumnWidth = 14.83

# 6) Update selection/scroll position left on the sheet by the author.
$ws.Range("G31").Select() | Out-Null
